$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values. A leading apostrophe forces each
# assignment to be stored as literal text (matching the source
# workbook, where Price/Volume/Coin/Link are all text cells) so
# numeric-looking strings like "1.00" or "7.42" are not silently
# coerced into Number cells, and without touching any cell styles.

$ws.Range('D2').Value = "'68.035.55"
$ws.Range('E2').Value = "'  +0.34%  "
$ws.Range('D3').Value = "'3.793.05"
$ws.Range('E3').Value = "'  -0.27%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "'  +0.21%  "
$ws.Range('D5').Value = "'602.10"
$ws.Range('E5').Value = "'  +0.94%  "
$ws.Range('D6').Value = "'165.36"
$ws.Range('E6').Value = "'  -1.06%  "
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('D8').Value = "'0.517"
$ws.Range('E8').Value = "'  -0.72%  "
$ws.Range('E9').Value = "'  -0.77%  "
$ws.Range('D10').Value = "'0.451"
$ws.Range('E10').Value = "'  +0.43%  "
$ws.Range('D11').Value = "'6.43"
$ws.Range('E11').Value = "'  +1.18%  "
$ws.Range('D12').Value = "'0.0000249"
$ws.Range('E12').Value = "'  -1.59%  "
$ws.Range('D13').Value = "'35.75"
$ws.Range('E13').Value = "'  -0.96%  "
$ws.Range('D14').Value = "'4.431.17"
$ws.Range('E14').Value = "'  -0.20%  "
$ws.Range('D15').Value = "'3.800.54"
$ws.Range('E15').Value = "'  -0.77%  "
$ws.Range('D16').Value = "'68.060.45"
$ws.Range('E16').Value = "'  +0.43%  "
$ws.Range('D17').Value = "'18.32"
$ws.Range('E17').Value = "'  -1.39%  "
$ws.Range('E18').Value = "'  +1.87%  "
$ws.Range('D19').Value = "'7.08"
$ws.Range('E19').Value = "'  -0.15%  "
$ws.Range('D20').Value = "'461.77"
$ws.Range('E20').Value = "'  +0.15%  "
$ws.Range('D21').Value = "'9.71"
$ws.Range('E21').Value = "'  -2.68%  "
$ws.Range('E22').Value = "'  -0.55%  "
$ws.Range('E23').Value = "'  -2.70%  "
$ws.Range('D24').Value = "'82.83"
$ws.Range('E24').Value = "'  -0.79%  "
$ws.Range('D25').Value = "'12.03"
$ws.Range('E25').Value = "'  -0.33%  "
$ws.Range('D26').Value = "'2.10"
$ws.Range('E26').Value = "'  +0.01%  "
$ws.Range('E27').Value = "'  -0.67%  "
$ws.Range('D28').Value = "'9.99"
$ws.Range('E28').Value = "'  -0.13%  "
$ws.Range('D29').Value = "'3.943.52"
$ws.Range('E29').Value = "'  -0.17%  "
$ws.Range('D30').Value = "'7.42"
$ws.Range('E30').Value = "'  +2.47%  "
$ws.Range('D31').Value = "'2.65"
$ws.Range('E31').Value = "'  -4.92%  "
$ws.Range('E32').Value = "'  -1.90%  "
$ws.Range('D33').Value = "'29.27"
$ws.Range('E33').Value = "'  -1.52%  "
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = "'  +0.15%  "
$ws.Range('D35').Value = "'8.99"
$ws.Range('E35').Value = "'  -0.78%  "
$ws.Range('D36').Value = "'0.0996"
$ws.Range('E36').Value = "'  -0.44%  "
$ws.Range('E37').Value = "'  +0.92%  "
$ws.Range('D38').Value = "'3.26"
$ws.Range('E38').Value = "'  -3.33%  "
$ws.Range('D39').Value = "'5.78"
$ws.Range('E39').Value = "'  +0.08%  "
$ws.Range('D40').Value = "'0.986"
$ws.Range('E40').Value = "'  -0.90%  "
$ws.Range('E41').Value = "'  +0.15%  "
$ws.Range('E42').Value = "'  +0.00%  "
$ws.Range('B43').Value = "'OKB"
$ws.Range('C43').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D43').Value = "'47.44"
$ws.Range('E43').Value = "'  -1.68%  "
$ws.Range('B44').Value = "'TheGraph"
$ws.Range('C44').Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range('D44').Value = "'0.298"
$ws.Range('E44').Value = "'  +0.30%  "
$ws.Range('D45').Value = "'43.12"
$ws.Range('E45').Value = "'  -1.47%  "
$ws.Range('D46').Value = "'151.01"
$ws.Range('E46').Value = "'  +0.29%  "
$ws.Range('D47').Value = "'8.34"
$ws.Range('E47').Value = "'  +0.46%  "
$ws.Range('E48').Value = "'  +2.07%  "
$ws.Range('D49').Value = "'395.40"
$ws.Range('E49').Value = "'  +1.08%  "
$ws.Range('D50').Value = "'27.05"
$ws.Range('E50').Value = "'  -0.80%  "
$ws.Range('D51').Value = "'1.34"
$ws.Range('E51').Value = "'  +5.61%  "
